# Convention change to support multi-axle vehicles:
#   "sAxleF" -> "sAxle1"      (row label, column A)
#   "Body_1Axle" -> "Body_Axle1"   (class label, column H)
# Applies to every sheet in the workbook (Trailer_Elula, Trailer_Elula_Unstable,
# Trailer_Thwala all share this layout).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("H4").Value = "Body_Axle1"
}
